$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking snapshot refresh (prices / 1h volume %, plus a couple of
# coins that moved rank and shifted rows 6-18 / 48-49 down one slot).
# Values are entered with a leading apostrophe so Excel stores them as
# literal text instead of auto-converting numeric-looking strings (like
# "312.02" or "1.98%") into Number/Percentage cells. The style is then
# reset to the sheet's untouched default (taken from F2, which every data
# row uses) so the quote-prefix flag that ".Value = ''...'" leaves behind
# does not show up as a spurious formatting change.
$defaultStyle = $ws.Range("F2").Style

$updates = [ordered]@{
    'D2' = '312.02'
    'E2' = '1.98%'
    'D3' = '37.33'
    'E3' = '0.43%'
    'D4' = '5.127'
    'E4' = '0.74%'
    'D5' = '0.07826'
    'E5' = '1.45%'
    'B6' = 'GateToken'
    'C6' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
    'D6' = '4.418'
    'E6' = '1.48%'
    'D7' = '8.277'
    'E7' = '0.93%'
    'B8' = 'FTXToken'
    'C8' = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
    'D8' = '1.902'
    'E8' = '-0.25%'
    'B9' = 'BTSEToken'
    'C9' = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
    'D9' = '2.836'
    'E9' = '-9.23%'
    'B10' = 'MXToken'
    'C10' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D10' = '0.9196'
    'E10' = '0.33%'
    'B11' = 'LiechtensteinCryptoassetsExchange'
    'C11' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'D11' = '0.1192'
    'E11' = '2.32%'
    'B12' = 'WazirX'
    'C12' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'D12' = '0.1918'
    'E12' = '2.81%'
    'B13' = 'MandalaExchangeToken'
    'C13' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'D13' = '0.08977'
    'E13' = '3.55%'
    'B14' = 'BitrueCoin'
    'C14' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'D14' = '0.03352'
    'E14' = '-1.54%'
    'B15' = 'BitMartToken'
    'C15' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'D15' = '0.09615'
    'E15' = '-0.74%'
    'B16' = 'BitForexToken'
    'C16' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'D16' = '0.001385'
    'E16' = '1.50%'
    'B17' = 'TigerCash'
    'C17' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'D17' = '0.005723'
    'E17' = '-3.34%'
    'B18' = 'LEO'
    'C18' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'D18' = '3.536'
    'E18' = '-2.05%'
    'D19' = '0.3442'
    'E19' = '0.94%'
    'D20' = '5.271'
    'E20' = '5.10%'
    'D21' = '0.1285'
    'E21' = '0.78%'
    'E22' = '0.04%'
    'D23' = '0.04366'
    'E23' = '0.94%'
    'D24' = '0.001252'
    'E24' = '3.15%'
    'D25' = '0.004668'
    'E25' = '2.91%'
    'D26' = '0.0001362'
    'E26' = '0.73%'
    'D27' = '0.0004003'
    'D39' = '0.02277'
    'E39' = '3.49%'
    'D40' = '0.05053'
    'E40' = '3.00%'
    'D41' = '0.007466'
    'E41' = '-0.90%'
    'D42' = '0.009088'
    'E42' = '-7.93%'
    'D43' = '0.1348'
    'E43' = '1.30%'
    'D44' = '0.001952'
    'E44' = '-2.19%'
    'D45' = '0.009320'
    'E45' = '9.77%'
    'D46' = '0.00006635'
    'E46' = '1.37%'
    'D47' = '0.00000000751'
    'E47' = '-0.01%'
    'B48' = 'BOLO'
    'C48' = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
    'D48' = '0.003278'
    'E48' = '9.33%'
    'B49' = 'CoinbaseStockToken'
    'C49' = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
    'D49' = '0.001003'
    'E49' = '-22.92%'
    'D50' = '0.00002103'
    'E50' = '-0.01%'
    'D51' = '0.0002002'
    'E51' = '-0.01%'
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.Value = "'" + $updates[$ref]
    $cell.Style = $defaultStyle
}

Write-Host "Applied $($updates.Count) cell updates"
